# Append a new row of statistics data to the "統計" (Statistics) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

# Find the next empty row after the existing data (row 12 -> new row 13)
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$ws.Cells.Item($newRow, 1).Value = "2025-08-28T18:25:28.704926"
$ws.Cells.Item($newRow, 2).Value = 13
$ws.Cells.Item($newRow, 3).Value = "全案件リスト"
$ws.Cells.Item($newRow, 4).Value = 53.8
$ws.Cells.Item($newRow, 5).Value = 6
$ws.Cells.Item($newRow, 6).Value = 4
$ws.Cells.Item($newRow, 7).Value = 13
